$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New committee member rows (51-63), matching the order the underlying
# source data was entered so that new shared strings are created in the
# same sequence as the target workbook.
$data = @(
  @("Björn", "Bornkamp", "Novartis", "CH"),
  @("Harlan", "Campbell", "U British Columbia", "CA"),
  @("David", "Dejardin", "Roche", "CH"),
  @("Jenny", "Devenport", "Roche", "CH"),
  @("Oliver", "Dukes", "U Ghent", "BE"),
  @("Elise", "Dumas", "EPFL", "CH"),
  @("Yair", "Goldberg", "Technion", "IL"),
  @("Dominic Edmund", "Magirr", "Novartis", "CH"),
  @("Francois", "Mercier", "Roche", "CH"),
  @("Tim P.", "Morris", "UCL", "UK"),
  @("Antonio", "Remiro-Azocar", "Novo Nordisk", "DK"),
  @("Garth", "Tarr", "U Sydney", "AU"),
  @("Marvin", "Wright", "BIPS", "DE")
)

$startRow = 51

# Step 1: row 51 is entered in full (first, last, institution, country)
$ws.Cells.Item($startRow, 1).Value = $data[0][0]
$ws.Cells.Item($startRow, 2).Value = $data[0][1]
$ws.Cells.Item($startRow, 3).Value = $data[0][2]
$ws.Cells.Item($startRow, 4).Value = $data[0][3]

# Step 2: rows 52-63 get their last name / institution / country first
for ($i = 1; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
    $ws.Cells.Item($r, 4).Value = $data[$i][3]
}

# Step 3: rows 52-63 then get their first names filled in
for ($i = 1; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
}

# Add new "sup_reviewer" column header in K1, copying the header style from A1
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "sup_reviewer"

# Flag Jack Kuipers (row 17) as a supporting reviewer
$ws.Range("K17").Value = 1

# Flag all the newly added committee members as supporting reviewers
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 11).Value = 1
}

$ws.Range("K2").Select()
